$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'53.513.70"
$ws.Range('D2').Style = $ws.Range('B2').Style
$ws.Range('E2').Value = '  +3.83%  '

# Row 3
$ws.Range('D3').Value = "'3.132.45"
$ws.Range('D3').Style = $ws.Range('B3').Style
$ws.Range('E3').Value = '  +2.28%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').Value = "'395.79"
$ws.Range('D5').Style = $ws.Range('B5').Style
$ws.Range('E5').Value = '  +2.58%  '

# Row 6
$ws.Range('D6').Value = "'108.61"
$ws.Range('D6').Style = $ws.Range('B6').Style
$ws.Range('E6').Value = '  +5.40%  '

# Row 7
$ws.Range('E7').Value = '  +0.16%  '

# Row 8
$ws.Range('E8').Value = '  -0.06%  '

# Row 9
$ws.Range('D9').Value = "'0.607"
$ws.Range('D9').Style = $ws.Range('B9').Style
$ws.Range('E9').Value = '  +3.85%  '

# Row 10
$ws.Range('D10').Value = "'38.47"
$ws.Range('D10').Style = $ws.Range('B10').Style
$ws.Range('E10').Value = '  +4.59%  '

# Row 11
$ws.Range('E11').Value = '  +1.22%  '

# Row 12
$ws.Range('E12').Value = '  +1.02%  '

# Row 13
$ws.Range('D13').Value = "'3.635.06"
$ws.Range('D13').Style = $ws.Range('B13').Style
$ws.Range('E13').Value = '  +2.26%  '

# Row 14
$ws.Range('D14').Value = "'18.91"
$ws.Range('D14').Style = $ws.Range('B14').Style
$ws.Range('E14').Value = '  +1.42%  '

# Row 15
$ws.Range('D15').Value = "'7.95"
$ws.Range('D15').Style = $ws.Range('B15').Style
$ws.Range('E15').Value = '  +2.21%  '

# Row 16
$ws.Range('E16').Value = '  +7.36%  '

# Row 17
$ws.Range('D17').Value = "'3.140.70"
$ws.Range('D17').Style = $ws.Range('B17').Style
$ws.Range('E17').Value = '  +2.26%  '

# Row 18
$ws.Range('D18').Value = "'10.43"
$ws.Range('D18').Style = $ws.Range('B18').Style
$ws.Range('E18').Value = '  -2.57%  '

# Row 19
$ws.Range('D19').Value = "'53.461.12"
$ws.Range('D19').Style = $ws.Range('B19').Style
$ws.Range('E19').Value = '  +3.56%  '

# Row 20
$ws.Range('D20').Value = "'3.24"
$ws.Range('D20').Style = $ws.Range('B20').Style
$ws.Range('E20').Value = '  +2.77%  '

# Row 21
$ws.Range('D21').Value = "'12.69"
$ws.Range('D21').Style = $ws.Range('B21').Style
$ws.Range('E21').Value = '  +1.82%  '

# Row 22
$ws.Range('D22').Value = "'0.0₃0970"
$ws.Range('D22').Style = $ws.Range('B22').Style
$ws.Range('E22').Value = '  +0.12%  '

# Row 23
$ws.Range('D23').Value = "'70.67"
$ws.Range('D23').Style = $ws.Range('B23').Style
$ws.Range('E23').Value = '  +0.71%  '

# Row 24
$ws.Range('D24').Value = "'269.79"
$ws.Range('D24').Style = $ws.Range('B24').Style
$ws.Range('E24').Value = '  +0.71%  '

# Row 25
$ws.Range('D25').Value = "'3.22"
$ws.Range('D25').Style = $ws.Range('B25').Style
$ws.Range('E25').Value = '  +2.45%  '

# Row 26
$ws.Range('D26').Value = "'7.88"
$ws.Range('D26').Style = $ws.Range('B26').Style
$ws.Range('E26').Value = '  -3.71%  '

# Row 27
$ws.Range('D27').Value = "'27.34"
$ws.Range('D27').Style = $ws.Range('B27').Style
$ws.Range('E27').Value = '  +1.63%  '

# Row 28
$ws.Range('D28').Value = "'7.32"
$ws.Range('D28').Style = $ws.Range('B28').Style
$ws.Range('E28').Value = '  +0.26%  '

# Row 29
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').Value = "'0.999"
$ws.Range('D29').Style = $ws.Range('B29').Style
$ws.Range('E29').Value = '  -0.16%  '

# Row 30
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = "'0.168"
$ws.Range('D30').Style = $ws.Range('B30').Style
$ws.Range('E30').Value = '  -1.35%  '

# Row 31
$ws.Range('E31').Value = '  +2.53%  '

# Row 32
$ws.Range('D32').Value = "'10.91"
$ws.Range('D32').Style = $ws.Range('B32').Style
$ws.Range('E32').Value = '  +6.42%  '

# Row 33
$ws.Range('E33').Value = '  +11.14%  '

# Row 34
$ws.Range('D34').Value = "'36.58"
$ws.Range('D34').Style = $ws.Range('B34').Style
$ws.Range('E34').Value = '  +5.07%  '

# Row 35
$ws.Range('D35').Value = "'2.07"
$ws.Range('D35').Style = $ws.Range('B35').Style
$ws.Range('E35').Value = '  +0.12%  '

# Row 36
$ws.Range('D36').Value = "'50.27"
$ws.Range('D36').Style = $ws.Range('B36').Style
$ws.Range('E36').Value = '  +0.46%  '

# Row 37
$ws.Range('D37').Value = "'3.62"
$ws.Range('D37').Style = $ws.Range('B37').Style
$ws.Range('E37').Value = '  +9.05%  '

# Row 38
$ws.Range('D38').Value = "'1.00"
$ws.Range('D38').Style = $ws.Range('B38').Style
$ws.Range('E38').Value = '  -0.08%  '

# Row 39
$ws.Range('D39').Value = "'2.76"
$ws.Range('D39').Style = $ws.Range('B39').Style
$ws.Range('E39').Value = '  +7.83%  '

# Row 40
$ws.Range('D40').Value = "'4.06"
$ws.Range('D40').Style = $ws.Range('B40').Style
$ws.Range('E40').Value = '  +8.58%  '

# Row 41
$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').Value = "'17.11"
$ws.Range('D41').Style = $ws.Range('B41').Style
$ws.Range('E41').Value = '  +1.31%  '

# Row 42
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').Value = "'0.286"
$ws.Range('D42').Style = $ws.Range('B42').Style
$ws.Range('E42').Value = '  -2.18%  '

# Row 43
$ws.Range('D43').Value = "'1.88"
$ws.Range('D43').Style = $ws.Range('B43').Style
$ws.Range('E43').Value = '  +0.47%  '

# Row 44
$ws.Range('D44').Value = "'129.93"
$ws.Range('D44').Style = $ws.Range('B44').Style
$ws.Range('E44').Value = '  +3.61%  '

# Row 45
$ws.Range('E45').Value = '  +0.88%  '

# Row 46
$ws.Range('D46').Value = "'21.97"
$ws.Range('D46').Style = $ws.Range('B46').Style
$ws.Range('E46').Value = '  +0.26%  '

# Row 47
$ws.Range('E47').Value = '  -0.22%  '

# Row 48
$ws.Range('E48').Value = '  -0.84%  '

# Row 49
$ws.Range('D49').Value = "'2.065.15"
$ws.Range('D49').Style = $ws.Range('B49').Style
$ws.Range('E49').Value = '  +1.41%  '

# Row 50
$ws.Range('D50').Value = "'0.0335"
$ws.Range('D50').Style = $ws.Range('B50').Style
$ws.Range('E50').Value = '  +4.76%  '

# Row 51
$ws.Range('D51').Value = "'0.0496"
$ws.Range('D51').Style = $ws.Range('B51').Style
$ws.Range('E51').Value = '  +16.28%  '
